$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 106 ---
# Column A: date serial changed from 45489.560787037 to 45489.2916666667
$ws.Cells.Item(106, 1).Value = 45489.2916666667

# --- Append new row 107 ---
$ws.Cells.Item(107, 1).Value = 45490.3697916667
$ws.Cells.Item(107, 2).Value = 300
$ws.Cells.Item(107, 3).Value = 6.11999988555908
$ws.Cells.Item(107, 4).Value = 6.11999988555908
$ws.Cells.Item(107, 5).Value = 6.11999988555908
$ws.Cells.Item(107, 6).Value = 6.11999988555908

# Column G (adj_close) is stored as text in this workbook (matches shared string "6.11999988555908").
# Prefix with an apostrophe so Excel stores it as text rather than a number, then reset
# the cell style back to Normal so no extra formatting/style gets attached to the cell.
$ws.Cells.Item(107, 7).Value = "'6.11999988555908"
$ws.Cells.Item(107, 7).Style = "Normal"

$ws.Cells.Item(107, 8).Value = "PAL.MI"

# Give the new date cell (A107) the same date/time display format as the rest of column A (e.g. A106).
$ws.Cells.Item(106, 1).Copy()
$ws.Cells.Item(107, 1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
